# Network visualization cleanup: remove duplicate / stray citation-network
# edges from the env_network sheet (rows that were accidentally duplicated
# when building the density/transitivity network).
#
# Original row layout (citingCase / citedCase) around the tail of the sheet,
# before this edit:
#   144: US - Tuna (Mexico) | Thailand - Cigarettes
#   145: US - Tuna (Mexico) | US - Section 337 Tariff Act   <- duplicate, remove
#   146: US - Tuna (Mexico) | EEC - Parts and Components    <- duplicate, remove
#   147: US - Tuna (Mexico) | Canada - FIRA                 <- duplicate, remove
#   148: US - Tuna (Mexico) | Canada - Herring and Salmon
#   149: US - Tuna (EEC)    | Canada - Herring and Salmon
#   150: US - Tuna (EEC)    | US - Canadian Tuna
#   151: US - Tuna (EEC)    | Canada - FIRA
#   152: US - Tuna (EEC)    | US - Section 337 Tariff Act
#   153: US - Tuna (EEC)    | US - Section 337 Tariff Act   <- duplicate, remove
#   154: US - Tuna (EEC)    | Thailand - Cigarettes
#   155: US - Taxes on Automobiles | Canada - Herring and Salmon
#
# Deleting rows 145, 146, 147 and 153 (working bottom-up so earlier deletes
# don't renumber the rows still queued for removal) reproduces the target
# state where the sheet ends at row 151.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(153, 147, 146, 145)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# Leave the selection where the author's cursor ended up after trimming the
# table.
[void]$ws.Range("B153").Select()
